$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (sp500 row), columns C..O get new correlation values
$ws.Range("C2").Value = -0.4739381494412931
$ws.Range("D2").Value = 0.1382547059950641
$ws.Range("E2").Value = 0.09155845137439581
$ws.Range("G2").Value = 0.03335360499461792
$ws.Range("H2").Value = 0.1328582935888218
$ws.Range("I2").Value = -0.07028726428208364
$ws.Range("J2").Value = 0.06559350788852032
$ws.Range("K2").Value = 0.07359892647541673
$ws.Range("L2").Value = -0.06826462012334182
$ws.Range("M2").Value = 0.03826994240463905
$ws.Range("N2").Value = 0.05534716907154695
$ws.Range("O2").Value = -0.008456246078771197

# Column B (sp500 column), mirrored symmetric values
$ws.Range("B3").Value = -0.4739381494412931
$ws.Range("B4").Value = 0.1382547059950641
$ws.Range("B5").Value = 0.09155845137439581
$ws.Range("B7").Value = 0.03335360499461792
$ws.Range("B8").Value = 0.1328582935888218
$ws.Range("B9").Value = -0.07028726428208364
$ws.Range("B10").Value = 0.06559350788852032
$ws.Range("B11").Value = 0.07359892647541673
$ws.Range("B12").Value = -0.06826462012334182
$ws.Range("B13").Value = 0.03826994240463905
$ws.Range("B14").Value = 0.05534716907154695
$ws.Range("B15").Value = -0.008456246078771197
